# Daily attendance processing - 2025-12-23 23:53:17
# Normalize the "Recorded By" column (G): when the value begins with the
# literal prefix "System, " move "System" from the front of the
# comma-separated list to the end of the list instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($null -ne $val -and $val.StartsWith("System, ")) {
        $rest = $val.Substring(8)
        $cell.Value = "$rest, System"
    }
}
